# ERD for Groupomania.xlsx -- add "Sheet2" (Customers / Products / Shops /
# Vendors / Sales mini-ERD tables), make it the active/selected sheet, and
# move Sheet1's selection from I7 to I8.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: move the selection cursor from I7 to I8 ---------------------
$ws1.Range("I8").Select()

# --- Add the new worksheet right after Sheet1 -----------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- Customers table (columns A-B) ---------------------------------------
$ws2.Range("A1").Value = "Customers"
$ws2.Range("A2").Value = "phone nr"
$ws2.Range("A3").Value = "customer nr"
$ws2.Range("A4").Value = "name"
$ws2.Range("A5").Value = "address"

# --- Products table (columns D-E) -----------------------------------------
$ws2.Range("D1").Value = "Products"
$ws2.Range("D2").Value = "price"
$ws2.Range("D3").Value = "type"
$ws2.Range("D4").Value = "manufacturer"

# --- Shops table (columns G-H) --------------------------------------------
$ws2.Range("G1").Value = "Shops"
$ws2.Range("G2").Value = "address"
$ws2.Range("G3").Value = "name"

# --- Vendors table (columns J-K) ------------------------------------------
$ws2.Range("J1").Value = "Vendors"
$ws2.Range("J2").Value = "staff number"
$ws2.Range("J3").Value = "name"

# --- Sales table (columns M-N) --------------------------------------------
$ws2.Range("M1").Value = "Sales"
$ws2.Range("M2").Value = "products"
$ws2.Range("M3").Value = "date"
$ws2.Range("M4").Value = "sum total"

# --- Highlight each table's title cell with a yellow fill -----------------
$ws2.Range("A1").Interior.Color = 65535
$ws2.Range("D1").Interior.Color = 65535
$ws2.Range("G1").Interior.Color = 65535
$ws2.Range("J1").Interior.Color = 65535
$ws2.Range("M1").Interior.Color = 65535

# --- Sheet2 is the tab that ends up selected/active -----------------------
$ws2.Activate()
$ws2.Range("A1").Select()
